# "separate nb for laptop" - update model3 results with values from a
# separate (laptop) notebook run, and refresh the per-column color-scale
# shading on the RMSE (D) and U (E) columns to match the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RgbColor($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# New R^2 / RMSE / U values per row (columns C, D, E), rows 2-9.
$data = @(
    @{ Row = 2;  C = -2.8566; D = 0.5356; E = 1.6822 },
    @{ Row = 3;  C = -1.1908; D = 0.6146; E = 1.4959 },
    @{ Row = 4;  C = -0.4115; D = 0.6587; E = 1.608  },
    @{ Row = 5;  C = -0.3611; D = 0.8083; E = 2.0039 },
    @{ Row = 6;  C = -0.4466; D = 0.9319; E = 2.3532 },
    @{ Row = 7;  C = -0.5543; D = 0.9698; E = 2.4605 },
    @{ Row = 8;  C = -0.6882; D = 1.0276; E = 2.5756 },
    @{ Row = 9;  C = -0.5314; D = 0.9981; E = 2.4839 }
)

# Color-scale shading (fill + font color) for the D (RMSE) and U (E) columns,
# keyed by cell address, matching the refreshed color scale for the new data.
$colors = @{
    "D2" = @{ Fill = (Set-RgbColor 0   68  27);  Font = (Set-RgbColor 241 241 241) };
    "D3" = @{ Fill = (Set-RgbColor 10  118 51);  Font = (Set-RgbColor 241 241 241) };
    "D4" = @{ Fill = (Set-RgbColor 35  139 69);  Font = (Set-RgbColor 0   0   0)   };
    "D5" = @{ Fill = (Set-RgbColor 135 205 134); Font = (Set-RgbColor 0   0   0)   };
    "D6" = @{ Fill = (Set-RgbColor 213 239 207); Font = (Set-RgbColor 0   0   0)   };
    "D7" = @{ Fill = (Set-RgbColor 230 245 225); Font = (Set-RgbColor 0   0   0)   };
    "D8" = @{ Fill = (Set-RgbColor 247 252 245); Font = (Set-RgbColor 0   0   0)   };
    "D9" = @{ Fill = (Set-RgbColor 239 249 235); Font = (Set-RgbColor 0   0   0)   };

    "E2" = @{ Fill = (Set-RgbColor 13  120 54);  Font = (Set-RgbColor 241 241 241) };
    "E3" = @{ Fill = (Set-RgbColor 0   68  27);  Font = (Set-RgbColor 241 241 241) };
    "E4" = @{ Fill = (Set-RgbColor 0   101 41);  Font = (Set-RgbColor 241 241 241) };
    "E5" = @{ Fill = (Set-RgbColor 104 190 112); Font = (Set-RgbColor 0   0   0)   };
    "E6" = @{ Fill = (Set-RgbColor 210 237 204); Font = (Set-RgbColor 0   0   0)   };
    "E7" = @{ Fill = (Set-RgbColor 232 246 227); Font = (Set-RgbColor 0   0   0)   };
    "E8" = @{ Fill = (Set-RgbColor 247 252 245); Font = (Set-RgbColor 0   0   0)   };
    "E9" = @{ Fill = (Set-RgbColor 235 247 231); Font = (Set-RgbColor 0   0   0)   };
}

foreach ($entry in $data) {
    $row = $entry.Row
    $ws.Range("C$row").Value = $entry.C
    $ws.Range("D$row").Value = $entry.D
    $ws.Range("E$row").Value = $entry.E
}

foreach ($addr in $colors.Keys) {
    $cell = $ws.Range($addr)
    $cell.Interior.Color = $colors[$addr].Fill
    $cell.Font.Color = $colors[$addr].Font
}
